# "Generate Report for Archive"
#
# The status "Ready for handoff" has moved on to "In Translation" for both
# locales, and the now-shorter status text means the Status column on each
# sheet re-sizes (narrower) to fit the new content.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E & F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
